$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on price cells whose new values look numeric,
# so Excel keeps storing them as text (matches existing column formatting).
$textCells = @("D4", "D5", "D6", "D7", "D8", "D10", "D15", "D18", "D19", "D23", "D24", "D28", "D31", "D33", "D34", "D36", "D38", "D39", "D41", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "63.923.65"
$ws.Range("E2").Value = "  -0.49%  "

# Row 3
$ws.Range("D3").Value = "3.127.55"
$ws.Range("E3").Value = "  -1.35%  "

# Row 4
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").Value = "566.70"
$ws.Range("E5").Value = "  -0.33%  "

# Row 6
$ws.Range("D6").Value = "160.59"
$ws.Range("E6").Value = "  -4.92%  "

# Row 7
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.18%  "

# Row 8
$ws.Range("D8").Value = "0.563"
$ws.Range("E8").Value = "  -7.73%  "

# Row 9
$ws.Range("E9").Value = "  -4.28%  "

# Row 10
$ws.Range("D10").Value = "6.53"
$ws.Range("E10").Value = "  -2.67%  "

# Row 11
$ws.Range("E11").Value = "  -1.84%  "

# Row 12
$ws.Range("D12").Value = "3.671.31"
$ws.Range("E12").Value = "  -1.33%  "

# Row 14
$ws.Range("D14").Value = "63.985.79"
$ws.Range("E14").Value = "  -0.46%  "

# Row 15
$ws.Range("D15").Value = "24.76"
$ws.Range("E15").Value = "  -2.52%  "

# Row 16
$ws.Range("D16").Value = "3.126.32"
$ws.Range("E16").Value = "  -1.64%  "

# Row 17
$ws.Range("E17").Value = "  -3.01%  "

# Row 18
$ws.Range("D18").Value = "399.65"
$ws.Range("E18").Value = "  -4.61%  "

# Row 19
$ws.Range("D19").Value = "12.47"
$ws.Range("E19").Value = "  -2.49%  "

# Row 20
$ws.Range("E20").Value = "  -2.88%  "

# Row 21
$ws.Range("E21").Value = "  +0.51%  "

# Row 22
$ws.Range("E22").Value = "  +3.37%  "

# Row 23
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.07%  "

# Row 24
$ws.Range("D24").Value = "67.81"
$ws.Range("E24").Value = "  -3.37%  "

# Row 25
$ws.Range("E25").Value = "  -1.90%  "

# Row 26
$ws.Range("E26").Value = "  -4.75%  "

# Row 27
$ws.Range("D27").Value = "0.0₂01000"
$ws.Range("E27").Value = "  -5.25%  "

# Row 28
$ws.Range("D28").Value = "8.69"
$ws.Range("E28").Value = "  -0.16%  "

# Row 29
$ws.Range("E29").Value = "  -0.05%  "

# Row 30
$ws.Range("E30").Value = "  -1.74%  "

# Row 31
$ws.Range("D31").Value = "21.00"
$ws.Range("E31").Value = "  -3.64%  "

# Row 32
$ws.Range("E32").Value = "  -1.90%  "

# Row 33
$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").Value = "157.42"
$ws.Range("E33").Value = "  +0.58%  "

# Row 34
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").Value = "4.77"
$ws.Range("E34").Value = "  -4.60%  "

# Row 35
$ws.Range("E35").Value = "  -3.10%  "

# Row 36
$ws.Range("D36").Value = "1.32"
$ws.Range("E36").Value = "  -3.19%  "

# Row 37
$ws.Range("D37").Value = "2.643.39"
$ws.Range("E37").Value = "  -2.07%  "

# Row 38
$ws.Range("D38").Value = "1.65"
$ws.Range("E38").Value = "  -2.26%  "

# Row 39
$ws.Range("D39").Value = "23.38"
$ws.Range("E39").Value = "  -4.27%  "

# Row 40
$ws.Range("E40").Value = "  -2.99%  "

# Row 41
$ws.Range("D41").Value = "0.686"
$ws.Range("E41").Value = "  -3.01%  "

# Row 42
$ws.Range("E42").Value = "  -2.04%  "

# Row 43
$ws.Range("D43").Value = "5.43"
$ws.Range("E43").Value = "  -4.88%  "

# Row 44
$ws.Range("D44").Value = "0.0252"
$ws.Range("E44").Value = "  -3.29%  "

# Row 45
$ws.Range("D45").Value = "284.69"
$ws.Range("E45").Value = "  -2.94%  "

# Row 46
$ws.Range("D46").Value = "20.86"
$ws.Range("E46").Value = "  -3.89%  "

# Row 47
$ws.Range("D47").Value = "0.997"
$ws.Range("E47").Value = "  -0.29%  "

# Row 48
$ws.Range("D48").Value = "0.0970"
$ws.Range("E48").Value = "  -2.28%  "

# Row 49
$ws.Range("D49").Value = "10.44"
$ws.Range("E49").Value = "  -0.04%  "

# Row 50
$ws.Range("D50").Value = "1.87"
$ws.Range("E50").Value = "  -7.11%  "

# Row 51
$ws.Range("D51").Value = "5.61"
$ws.Range("E51").Value = "  -2.61%  "
